$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# --- Row 15/16: swap ShibaInu and WrappedliquidstakedEther2.0 rows ---
Set-TextValue "B15" "WrappedliquidstakedEther2.0"
Set-TextValue "C15" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D15" "2.811.33"
Set-TextValue "E15" "  +0.30%  "
Set-TextValue "B16" "ShibaInu"
Set-TextValue "C16" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D16" "0.0000167"
Set-TextValue "E16" "  -1.63%  "

# --- Remaining price / volume(1h) updates ---
Set-TextValue "D2" "60.203.95"
Set-TextValue "E2" "  -0.91%  "
Set-TextValue "D3" "2.380.15"
Set-TextValue "E3" "  -0.85%  "
Set-TextValue "D5" "561.37"
Set-TextValue "D6" "138.79"
Set-TextValue "E6" "  -0.37%  "
Set-TextValue "E8" "  +0.59%  "
Set-TextValue "D9" "2.379.16"
Set-TextValue "E9" "  -0.06%  "
Set-TextValue "E10" "  -1.22%  "
Set-TextValue "E12" "  +0.78%  "
Set-TextValue "E13" "  +0.04%  "
Set-TextValue "D14" "25.79"
Set-TextValue "E14" "  -0.28%  "
Set-TextValue "D17" "59.988.32"
Set-TextValue "E17" "  -1.26%  "
Set-TextValue "D18" "2.371.05"
Set-TextValue "E18" "  -0.59%  "
Set-TextValue "D19" "8.05"
Set-TextValue "E19" "  +11.73%  "
Set-TextValue "E20" "  +0.06%  "
Set-TextValue "D21" "322.14"
Set-TextValue "E21" "  +0.25%  "
Set-TextValue "E22" "  +1.25%  "
Set-TextValue "D23" "6.05"
Set-TextValue "E23" "  -0.31%  "
Set-TextValue "E24" "  -0.02%  "
Set-TextValue "E25" "  -1.27%  "
Set-TextValue "D26" "64.08"
Set-TextValue "E26" "  -0.35%  "
Set-TextValue "D27" "560.67"
Set-TextValue "E27" "  -1.79%  "
Set-TextValue "D28" "8.14"
Set-TextValue "E28" "  -4.61%  "
Set-TextValue "D29" "2.498.10"
Set-TextValue "E29" "  -0.25%  "
Set-TextValue "D30" "0.0₃0930"
Set-TextValue "E30" "  +2.45%  "
Set-TextValue "D31" "7.99"
Set-TextValue "E31" "  +2.21%  "
Set-TextValue "E32" "  -2.10%  "
Set-TextValue "E33" "  -2.00%  "
Set-TextValue "E34" "  +1.02%  "
Set-TextValue "E35" "  -0.57%  "
Set-TextValue "E36" "  +5.19%  "
Set-TextValue "D37" "153.67"
Set-TextValue "E37" "  +3.93%  "
Set-TextValue "E38" "  +0.05%  "
Set-TextValue "E39" "  -0.77%  "
Set-TextValue "E40" "  +0.09%  "
Set-TextValue "D41" "5.06"
Set-TextValue "E41" "  -0.12%  "
Set-TextValue "E42" "  -0.12%  "
Set-TextValue "D43" "41.60"
Set-TextValue "E43" "  +0.33%  "
Set-TextValue "D44" "1.65"
Set-TextValue "E44" "  -0.22%  "
Set-TextValue "E45" "  +4.95%  "
Set-TextValue "D46" "0.0₆0287"
Set-TextValue "E46" "  +0.62%  "
Set-TextValue "D47" "140.68"
Set-TextValue "E47" "  +0.41%  "
Set-TextValue "E48" "  +0.94%  "
Set-TextValue "D49" "0.586"
Set-TextValue "E49" "  +0.43%  "
Set-TextValue "D50" "0.0501"
Set-TextValue "E50" "  -0.23%  "
Set-TextValue "D51" "19.23"
Set-TextValue "E51" "  -0.07%  "
